$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Status text: handback is complete, so the "Ready for handoff" status
# becomes "Handed back: in sync with en-US" everywhere it is shown.
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value     = "Handed back: in sync with en-US"
$dede.Range("C2").Value     = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: record the generated handback package + timestamp.
# ---------------------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e137cb76b1c0ff7456170f8c7a5e471d38acfdb/e2e/7a76190d-de93-40ea-8880-55568815d466.md", "", "", "7a76190d-de93-40ea-8880-55568815d466.md") | Out-Null
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("J2").Value = "7a76190d-de93-40ea-8880-55568815d466.09e8e78fff570da572b61e0f464184ed810dbd09.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 07:07:44"

# ---------------------------------------------------------------------------
# de-de sheet: record the generated handback package + timestamp.
# ---------------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e137cb76b1c0ff7456170f8c7a5e471d38acfdb/e2e/7a76190d-de93-40ea-8880-55568815d466.md", "", "", "7a76190d-de93-40ea-8880-55568815d466.md") | Out-Null
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Range("J2").Value = "7a76190d-de93-40ea-8880-55568815d466.09e8e78fff570da572b61e0f464184ed810dbd09.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 07:07:52"

# ---------------------------------------------------------------------------
# Widen the "Status"/"Latest Target File"/"Latest Handback File" columns so
# the new, longer status/file-name text is fully visible.
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

$zhcn.Columns.Item(3).ColumnWidth  = 29.1
$zhcn.Columns.Item(9).ColumnWidth  = 39.1
$zhcn.Columns.Item(10).ColumnWidth = 39.1

$dede.Columns.Item(3).ColumnWidth  = 29.1
$dede.Columns.Item(9).ColumnWidth  = 39.1
$dede.Columns.Item(10).ColumnWidth = 39.1
